# Verify_40V_Calculation_For_FIM.xlsx - "Updated test data as per new implemenation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string "40V (A)" (used by H8 and H9) becomes "40V Rail(A)"
$ws.Range("H8").Value = "40V Rail(A)"
$ws.Range("H9").Value = "40V Rail(A)"

# Leave the selection where the author ended up when they saved the file
$ws.Range("H9").Select() | Out-Null
